# Saccade extraction now working and displaying.
#
# Remove four now-obsolete "to do" paragraphs from the notes:
#   - "Need to clip out average traces centered on peak velocity and aligned"
#   - "From this position, a maximum of +-4 can be achieved."
#   - "Do blocks of 4 (1:4) selecting amplitudes in random order"
#   - "Ignore the sign of the amplitudes, so it can effectively be used to keep centered"
#
# Delete from the bottom of the document upward so earlier paragraph
# indices stay valid as later ones are removed.

$d = $word.ActiveDocument

$targets = @(
    "From this position, a maximum of +-4 can be achieved.",
    "Do blocks of 4 (1:4) selecting amplitudes in random order",
    "Ignore the sign of the amplitudes, so it can effectively be used to keep centered",
    "Need to clip out average traces centered on peak velocity and aligned"
)

foreach ($target in $targets) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}
